$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.994.45'
$ws.Cells.Item(2, 5).Value = '  -0.37%  '
$ws.Cells.Item(3, 4).Value = '1.859.63'
$ws.Cells.Item(3, 5).Value = '  -0.87%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '311.97'
$ws.Cells.Item(5, 5).Value = '  -0.47%  '
$ws.Cells.Item(6, 5).Value = '  +0.11%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5150'
$ws.Cells.Item(7, 5).Value = '  +1.35%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3836'
$ws.Cells.Item(8, 5).Value = '  -0.27%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08239'
$ws.Cells.Item(9, 5).Value = '  -9.94%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.111'
$ws.Cells.Item(10, 5).Value = '  -1.05%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '41.49'
$ws.Cells.Item(11, 5).Value = '  -0.18%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.193'
$ws.Cells.Item(12, 5).Value = '  -2.54%  '
$ws.Cells.Item(13, 5).Value = '  -1.08%  '
$ws.Cells.Item(14, 4).Value = '1.865.90'
$ws.Cells.Item(14, 5).Value = '  -0.22%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.256'
$ws.Cells.Item(15, 5).Value = '  +0.60%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.003'
$ws.Cells.Item(16, 5).Value = '  +0.04%  '
$ws.Cells.Item(17, 5).Value = '  -1.90%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '90.59'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06650'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.67'
$ws.Cells.Item(20, 5).Value = '  -2.75%  '
$ws.Cells.Item(21, 5).Value = '  +0.08%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.004'
$ws.Cells.Item(22, 5).Value = '  -1.80%  '
$ws.Cells.Item(23, 4).Value = '28.021.75'
$ws.Cells.Item(23, 5).Value = '  -0.35%  '
$ws.Cells.Item(24, 5).Value = '  -3.25%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.262'
$ws.Cells.Item(25, 5).Value = '  -0.68%  '
$ws.Cells.Item(26, 4).Value = '2.074.22'
$ws.Cells.Item(26, 5).Value = '  -0.44%  '
$ws.Cells.Item(27, 5).Value = '  -2.19%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '157.69'
$ws.Cells.Item(28, 5).Value = '  +0.07%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '20.46'
$ws.Cells.Item(29, 5).Value = '  -1.59%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '124.54'
$ws.Cells.Item(30, 5).Value = '  -1.81%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.1064'
$ws.Cells.Item(31, 5).Value = '  +1.26%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.029'
$ws.Cells.Item(32, 5).Value = '  -3.22%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.925'
$ws.Cells.Item(33, 5).Value = '  +5.61%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.595'
$ws.Cells.Item(34, 5).Value = '  -0.25%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.381'
$ws.Cells.Item(35, 5).Value = '  -3.29%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02416'
$ws.Cells.Item(36, 5).Value = '  -0.94%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06499'
$ws.Cells.Item(37, 5).Value = '  -1.27%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2177'
$ws.Cells.Item(38, 5).Value = '  +0.06%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.6549'
$ws.Cells.Item(39, 5).Value = '  +2.23%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.196'
$ws.Cells.Item(40, 5).Value = '  -1.19%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.985'
$ws.Cells.Item(41, 5).Value = '  +1.43%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.214'
$ws.Cells.Item(42, 5).Value = '  -1.89%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.16'
$ws.Cells.Item(43, 5).Value = '  -3.64%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6162'
$ws.Cells.Item(44, 5).Value = '  +2.40%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.06'
$ws.Cells.Item(45, 5).Value = '  -1.07%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.283'
$ws.Cells.Item(46, 5).Value = '  +0.45%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.668'
$ws.Cells.Item(47, 5).Value = '  -0.08%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.006'
$ws.Cells.Item(48, 5).Value = '  +0.31%  '
$ws.Cells.Item(49, 5).Value = '  -1.54%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '120.61'
$ws.Cells.Item(50, 5).Value = '  -0.67%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '78.26'
$ws.Cells.Item(51, 5).Value = '  -1.96%  '
